$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 74

# Columns A-D hold text values ("2025-02-17", "22:56:21", "Monday", "07").
# Force text formatting first so Excel doesn't auto-convert the date/time
# strings into date/time serials or drop the leading zero on "07".
$ws.Range("A" + $row + ":D" + $row).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-02-17"
$ws.Cells.Item($row, 2).Value = "22:56:21"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).Value = "07"
$ws.Cells.Item($row, 5).Value = 128974
$ws.Cells.Item($row, 6).Value = 140386
$ws.Cells.Item($row, 7).Value = 171086
$ws.Cells.Item($row, 8).Value = 159221
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 145418
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 192328
$ws.Cells.Item($row, 14).Value = 115161
$ws.Cells.Item($row, 15).Value = 45468
$ws.Cells.Item($row, 16).Value = 28886
$ws.Cells.Item($row, 17).Value = 66727
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 46753
$ws.Cells.Item($row, 20).Value = -1
